# Auto-generated script applying numeric cell updates per the commit diff.
# Sheets map: sheet1=ALC, sheet2=ARM, sheet3=BSM, sheet4=CRP, sheet5=CUL, sheet6=GSM, sheet7=LTW, sheet8=WVR

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 159
$ws.Range("I5").Value = 116.375
$ws.Range("K5").Value = 116.375
$ws.Range("M5").Value = -1.375
$ws.Range("H32").Value = 741.65
$ws.Range("I32").Value = 719.4286
$ws.Range("J32").Value = 753.61536
$ws.Range("K32").Value = 719.4286
$ws.Range("L32").Value = 753.61536
$ws.Range("M32").Value = -393.4286
$ws.Range("N32").Value = -1405.61536
$ws.Range("H33").Value = 1976369.2
$ws.Range("I33").Value = 2392435.8
$ws.Range("J33").Value = 54.25
$ws.Range("K33").Value = 2392435.8
$ws.Range("L33").Value = 54.25
$ws.Range("M33").Value = -2392206.8
$ws.Range("N33").Value = -512.25
$ws.Range("H112").Value = 1274.1212
$ws.Range("J112").Value = 1274.1212
$ws.Range("L112").Value = 3822.3636
$ws.Range("N112").Value = -6038.363600000001
$ws.Range("H125").Value = 1149.6666
$ws.Range("J125").Value = 1199.6
$ws.Range("L125").Value = 10796.4
$ws.Range("N125").Value = -15716.4
$ws.Range("H129").Value = 1461.7301
$ws.Range("J129").Value = 1606.1072
$ws.Range("L129").Value = 4818.321599999999
$ws.Range("N129").Value = -14818.3216
$ws.Range("H132").Value = 37186136
$ws.Range("I132").Value = 45634310
$ws.Range("K132").Value = 136902930
$ws.Range("M132").Value = -136900400
$ws.Range("H138").Value = 2275.1538
$ws.Range("I138").Value = 1152.7222
$ws.Range("J138").Value = 3237.238
$ws.Range("K138").Value = 3458.1666
$ws.Range("L138").Value = 9711.714
$ws.Range("M138").Value = 1681.8334
$ws.Range("N138").Value = -19991.714

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3682.1785
$ws.Range("I32").Value = 3802.3386
$ws.Range("K32").Value = 3802.3386
$ws.Range("M32").Value = -3515.3386
$ws.Range("H61").Value = 1738.2
$ws.Range("I61").Value = 1738.2
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1738.2
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1526.2
$ws.Range("N61").Value = $null
$ws.Range("H74").Value = 381663.3
$ws.Range("I74").Value = 592400.9
$ws.Range("J74").Value = 2335.7
$ws.Range("K74").Value = 592400.9
$ws.Range("L74").Value = 2335.7
$ws.Range("M74").Value = -591526.9
$ws.Range("N74").Value = -4083.7
$ws.Range("H77").Value = 381663.3
$ws.Range("I77").Value = 592400.9
$ws.Range("J77").Value = 2335.7
$ws.Range("K77").Value = 2962004.5
$ws.Range("L77").Value = 11678.5
$ws.Range("M77").Value = -2957636.5
$ws.Range("N77").Value = -20414.5
$ws.Range("H122").Value = 2418.762
$ws.Range("I122").Value = 1372.9333
$ws.Range("J122").Value = 5033.3335
$ws.Range("K122").Value = 4118.7999
$ws.Range("L122").Value = 15100.0005
$ws.Range("M122").Value = -1668.7999
$ws.Range("N122").Value = -20000.0005
$ws.Range("H132").Value = 2849.4773
$ws.Range("I132").Value = 1941.4706
$ws.Range("J132").Value = 3421.1853
$ws.Range("K132").Value = 5824.4118
$ws.Range("L132").Value = 10263.5559
$ws.Range("M132").Value = -3294.4118
$ws.Range("N132").Value = -15323.5559
$ws.Range("H136").Value = 1738.2
$ws.Range("I136").Value = 1738.2
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 5214.6
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -2664.6
$ws.Range("N136").Value = $null

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4461
$ws.Range("I134").Value = 1412
$ws.Range("J134").Value = 6102.769
$ws.Range("K134").Value = 4236
$ws.Range("L134").Value = 18308.307
$ws.Range("M134").Value = -1701
$ws.Range("N134").Value = -23378.307

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2475.9592
$ws.Range("I31").Value = 920.72
$ws.Range("K31").Value = 920.72
$ws.Range("M31").Value = -625.72
$ws.Range("H34").Value = 2475.9592
$ws.Range("I34").Value = 920.72
$ws.Range("K34").Value = 920.72
$ws.Range("M34").Value = -718.72

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 376.92307
$ws.Range("I18").Value = 81.818184
$ws.Range("J18").Value = 2000
$ws.Range("K18").Value = 245.454552
$ws.Range("L18").Value = 6000
$ws.Range("M18").Value = -76.45455200000001
$ws.Range("N18").Value = -6338
$ws.Range("H68").Value = 1210.9143
$ws.Range("I68").Value = 721.2174
$ws.Range("J68").Value = 2149.5
$ws.Range("K68").Value = 2163.6522
$ws.Range("L68").Value = 6448.5
$ws.Range("M68").Value = -1352.6522
$ws.Range("N68").Value = -8070.5
$ws.Range("H71").Value = 1210.9143
$ws.Range("I71").Value = 721.2174
$ws.Range("J71").Value = 2149.5
$ws.Range("K71").Value = 6490.9566
$ws.Range("L71").Value = 19345.5
$ws.Range("M71").Value = -2434.9566
$ws.Range("N71").Value = -27457.5
$ws.Range("H113").Value = 589.1795
$ws.Range("I113").Value = 490.4
$ws.Range("J113").Value = 650.9167
$ws.Range("K113").Value = 1471.2
$ws.Range("L113").Value = 1952.7501
$ws.Range("M113").Value = 698.8000000000002
$ws.Range("N113").Value = -6292.7501
$ws.Range("H131").Value = 765.16
$ws.Range("I131").Value = 476
$ws.Range("J131").Value = 793.75824
$ws.Range("K131").Value = 1428
$ws.Range("L131").Value = 2381.27472
$ws.Range("M131").Value = 3612
$ws.Range("N131").Value = -12461.27472

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5492.113
$ws.Range("I70").Value = 5247.2354
$ws.Range("K70").Value = 5247.2354
$ws.Range("M70").Value = -4977.2354
$ws.Range("H73").Value = 5492.113
$ws.Range("I73").Value = 5247.2354
$ws.Range("K73").Value = 5247.2354
$ws.Range("M73").Value = -4311.2354
$ws.Range("H126").Value = 3879.8
$ws.Range("I126").Value = 2888.611
$ws.Range("J126").Value = 6428.5713
$ws.Range("K126").Value = 8665.832999999999
$ws.Range("L126").Value = 19285.7139
$ws.Range("M126").Value = -6195.832999999999
$ws.Range("N126").Value = -24225.7139
$ws.Range("H132").Value = 4642.08
$ws.Range("I132").Value = 3799.6365
$ws.Range("J132").Value = 5304
$ws.Range("K132").Value = 11398.9095
$ws.Range("L132").Value = 15912
$ws.Range("M132").Value = -8868.9095
$ws.Range("N132").Value = -20972

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 59857.5
$ws.Range("J36").Value = 59857.5
$ws.Range("L36").Value = 59857.5
$ws.Range("N36").Value = -60981.5
$ws.Range("H82").Value = 5060.875
$ws.Range("I82").Value = 5588.619
$ws.Range("J82").Value = 1366.6666
$ws.Range("K82").Value = 5588.619
$ws.Range("L82").Value = 1366.6666
$ws.Range("M82").Value = -5227.619
$ws.Range("N82").Value = -2088.6666
$ws.Range("H85").Value = 5060.875
$ws.Range("I85").Value = 5588.619
$ws.Range("J85").Value = 1366.6666
$ws.Range("K85").Value = 5588.619
$ws.Range("L85").Value = 1366.6666
$ws.Range("M85").Value = -4340.619
$ws.Range("N85").Value = -3862.6666
$ws.Range("H100").Value = 5055.8887
$ws.Range("I100").Value = 1900.6
$ws.Range("J100").Value = 9000
$ws.Range("K100").Value = 1900.6
$ws.Range("L100").Value = 9000
$ws.Range("M100").Value = -1359.6
$ws.Range("N100").Value = -10082
$ws.Range("H122").Value = 3909.087
$ws.Range("I122").Value = 2800.8
$ws.Range("J122").Value = 4216.9443
$ws.Range("K122").Value = 8402.400000000001
$ws.Range("L122").Value = 12650.8329
$ws.Range("M122").Value = -5952.400000000001
$ws.Range("N122").Value = -17550.8329
$ws.Range("H132").Value = 6262
$ws.Range("I132").Value = 4013.25
$ws.Range("K132").Value = 12039.75
$ws.Range("M132").Value = -9509.75

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4712.2085
$ws.Range("I122").Value = 2474.75
$ws.Range("J122").Value = 5830.9375
$ws.Range("K122").Value = 7424.25
$ws.Range("L122").Value = 17492.8125
$ws.Range("M122").Value = -4974.25
$ws.Range("N122").Value = -22392.8125
$ws.Range("H132").Value = 7411688.5
$ws.Range("I132").Value = 5653
$ws.Range("J132").Value = 14495722
$ws.Range("K132").Value = 16959
$ws.Range("L132").Value = 43487166
$ws.Range("M132").Value = -14429
$ws.Range("N132").Value = -43492226
$ws.Range("H136").Value = 12995.866
$ws.Range("I136").Value = 16294
$ws.Range("J136").Value = 10110
$ws.Range("K136").Value = 48882
$ws.Range("L136").Value = 30330
$ws.Range("M136").Value = -46332
$ws.Range("N136").Value = -35430
